# Append new vocabulary rows (59-64) to the active sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @("infringement", "noun", "violação", 0),
    @("inextricable", "adjective", "inextricável", 0),
    @("outrage", "noun", "ultraje", 0),
    @("vitriol", "noun", "vitríolo", 0),
    @("towards", "preposition", "em relação", 0),
    @("refrain", "verb", "abster", 0)
)

$startRow = 59
for ($i = 0; $i -lt $data.Length; $i++) {
    $r = $startRow + $i
    $row = $data[$i]
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
}
